$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1143363072"
$ws.Range("D16").Value = "CINDY DAYANA MADRID ORTIZ"
$ws.Range("E16").Value = "1802"
$ws.Range("F16").Value = 80000
$ws.Range("G16").Value = 2000000

$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1143363072"
$ws.Range("D17").Value = "CINDY DAYANA MADRID ORTIZ"
$ws.Range("E17").Value = "1801"
$ws.Range("F17").Value = 80000
$ws.Range("G17").Value = 2000000

$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1129523066"
$ws.Range("D18").Value = "BRENDA SOFIA STRUX MONTERROZA"
$ws.Range("E18").Value = "1802"
$ws.Range("F18").Value = 48000
$ws.Range("G18").Value = 1200000

$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1129523066"
$ws.Range("D19").Value = "BRENDA SOFIA STRUX MONTERROZA"
$ws.Range("E19").Value = "1801"
$ws.Range("F19").Value = 48000
$ws.Range("G19").Value = 1200000

$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "1143332454"
$ws.Range("D20").Value = "PABLO RAFAEL HERRERA CAPDEVILLA"
$ws.Range("E20").Value = "1802"
$ws.Range("F20").Value = 100000
$ws.Range("G20").Value = 2500000

$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "1143332454"
$ws.Range("D21").Value = "PABLO RAFAEL HERRERA CAPDEVILLA"
$ws.Range("E21").Value = "1801"
$ws.Range("F21").Value = 100000
$ws.Range("G21").Value = 2500000

$ws.Range("B22").Value = "CC"
$ws.Range("C22").Value = "1047435092"
$ws.Range("D22").Value = "ANA ISABEL VELASCO BARRETO"
$ws.Range("E22").Value = "1802"
$ws.Range("F22").Value = 80000
$ws.Range("G22").Value = 2000000

$ws.Range("B23").Value = "CC"
$ws.Range("C23").Value = "1047435092"
$ws.Range("D23").Value = "ANA ISABEL VELASCO BARRETO"
$ws.Range("E23").Value = "1801"
$ws.Range("F23").Value = 80000
$ws.Range("G23").Value = 2000000

